# Insert a new weekly price-record row at row 30 of the "Arveja Verde" sheet.
# All rows currently at 30..101 shift down to 31..102 (handled automatically
# by Rows.Insert), and the new row 30 gets a fresh record (same market /
# product metadata as every other row, new date + price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Range("A30").Value2 = 3
$ws.Range("B30").Value2 = "Femacal de La Calera"
$ws.Range("C30").Value2 = "Coquimbo"
$ws.Range("D30").Value2 = 45246
$ws.Range("E30").Value2 = 5
$ws.Range("F30").Value2 = 100112022
$ws.Range("G30").Value2 = "Arveja Verde"
$ws.Range("H30").Value2 = "Perfection"
$ws.Range("I30").Value2 = "Primera"
$ws.Range("J30").Value2 = 35
$ws.Range("K30").Value2 = 35000
$ws.Range("L30").Value2 = 35000
$ws.Range("M30").Value2 = 35000
$ws.Range("N30").Value2 = "$/saco 25 kilos"
$ws.Range("O30").Value2 = "Provincia de Limarí"
$ws.Range("P30").Value2 = 1400
$ws.Range("Q30").Value2 = 25
$ws.Range("R30").Value2 = "Hortaliza"
